$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H116").Value = 5382.3335
$ws.Range("J116").Value = 4443
$ws.Range("L116").Value = 4443
$ws.Range("N116").Value = -11327
$ws.Range("H125").Value = 4039.25
$ws.Range("J125").Value = 4831.5
$ws.Range("L125").Value = 43483.5
$ws.Range("N125").Value = -48403.5
$ws.Range("H141").Value = 1218.8846
$ws.Range("I141").Value = 1179.64
$ws.Range("J141").Value = 2200
$ws.Range("K141").Value = 3538.92
$ws.Range("L141").Value = 6600
$ws.Range("M141").Value = 1641.08
$ws.Range("N141").Value = -16960

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 44974.5
$ws.Range("J55").Value = 49949.5
$ws.Range("L55").Value = 49949.5
$ws.Range("N55").Value = -50579.5
$ws.Range("H74").Value = 126838.875
$ws.Range("I74").Value = 126838.875
$ws.Range("K74").Value = 126838.875
$ws.Range("M74").Value = -125964.875
$ws.Range("H77").Value = 126838.875
$ws.Range("I77").Value = 126838.875
$ws.Range("K77").Value = 634194.375
$ws.Range("M77").Value = -629826.375
$ws.Range("H97").Value = 9917.429
$ws.Range("I97").Value = 14094.875
$ws.Range("K97").Value = 14094.875
$ws.Range("M97").Value = -13598.875
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 3249.6316
$ws.Range("I110").Value = 3185.75
$ws.Range("J110").Value = 3296.0908
$ws.Range("K110").Value = 3185.75
$ws.Range("L110").Value = 3296.0908
$ws.Range("M110").Value = -1140.75
$ws.Range("N110").Value = -7386.0908
$ws.Range("H122").Value = 2113.524
$ws.Range("I122").Value = 1699.3846
$ws.Range("K122").Value = 5098.1538
$ws.Range("M122").Value = -2648.1538
$ws.Range("N122").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2466
$ws.Range("I20").Value = 2182.842
$ws.Range("K20").Value = 2182.842
$ws.Range("M20").Value = -1935.842
$ws.Range("N20").ClearContents()
$ws.Range("H86").Value = 4590.727
$ws.Range("I86").Value = 4385.5713
$ws.Range("J86").Value = 4949.75
$ws.Range("K86").Value = 4385.5713
$ws.Range("L86").Value = 4949.75
$ws.Range("M86").Value = -3262.5713
$ws.Range("N86").Value = -7195.75
$ws.Range("H89").Value = 4590.727
$ws.Range("I89").Value = 4385.5713
$ws.Range("J89").Value = 4949.75
$ws.Range("K89").Value = 21927.8565
$ws.Range("L89").Value = 24748.75
$ws.Range("M89").Value = -16311.8565
$ws.Range("N89").Value = -35980.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1301
$ws.Range("I3").Value = 734.6667
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 734.6667
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -621.6667
$ws.Range("N3").Value = -3226
$ws.Range("H12").Value = 4311.8
$ws.Range("I12").Value = 4311.8
$ws.Range("K12").Value = 4311.8
$ws.Range("M12").Value = -4141.8
$ws.Range("H13").Value = 50
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H23").Value = 3999.5
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H27").Value = 3999.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H31").Value = 2615.7
$ws.Range("I31").Value = 1989.6875
$ws.Range("K31").Value = 1989.6875
$ws.Range("M31").Value = -1694.6875
$ws.Range("N31").ClearContents()
$ws.Range("H33").Value = 1819.75
$ws.Range("I33").Value = 1819.75
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1819.75
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1440.75
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 2615.7
$ws.Range("I34").Value = 1989.6875
$ws.Range("K34").Value = 1989.6875
$ws.Range("M34").Value = -1787.6875
$ws.Range("N34").ClearContents()
$ws.Range("H38").Value = 15308.4
$ws.Range("J38").Value = 15308.4
$ws.Range("L38").Value = 15308.4
$ws.Range("N38").Value = -16062.4
$ws.Range("H39").Value = 16599.5
$ws.Range("I39").Value = 8799.666999999999
$ws.Range("K39").Value = 8799.666999999999
$ws.Range("M39").Value = -8408.666999999999
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 15308.4
$ws.Range("J46").Value = 15308.4
$ws.Range("L46").Value = 15308.4
$ws.Range("N46").Value = -15730.4
$ws.Range("H49").Value = 16599.5
$ws.Range("I49").Value = 8799.666999999999
$ws.Range("K49").Value = 8799.666999999999
$ws.Range("M49").Value = -8617.666999999999
$ws.Range("N49").ClearContents()
$ws.Range("H132").Value = 1783.0454
$ws.Range("I132").Value = 1783.0454
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5349.1362
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2819.1362
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 137.44444
$ws.Range("I2").Value = 135.875
$ws.Range("K2").Value = 135.875
$ws.Range("M2").Value = -22.875
$ws.Range("N2").ClearContents()
$ws.Range("H80").Value = 4805.8335
$ws.Range("I80").Value = 5105.273
$ws.Range("J80").Value = 4335.2856
$ws.Range("K80").Value = 5105.273
$ws.Range("L80").Value = 4335.2856
$ws.Range("M80").Value = -4107.273
$ws.Range("N80").Value = -6331.2856
$ws.Range("H83").Value = 4805.8335
$ws.Range("I83").Value = 5105.273
$ws.Range("J83").Value = 4335.2856
$ws.Range("K83").Value = 25526.365
$ws.Range("L83").Value = 21676.428
$ws.Range("M83").Value = -20534.365
$ws.Range("N83").Value = -31660.428
$ws.Range("H122").Value = 4848.625
$ws.Range("I122").Value = 1756.4
$ws.Range("K122").Value = 5269.200000000001
$ws.Range("M122").Value = -2819.200000000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 8013.5835
$ws.Range("I126").Value = 8169.2856
$ws.Range("K126").Value = 24507.8568
$ws.Range("M126").Value = -22037.8568
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11104.056
$ws.Range("I7").Value = 15200.091
$ws.Range("J7").Value = 4667.4287
$ws.Range("K7").Value = 15200.091
$ws.Range("L7").Value = 4667.4287
$ws.Range("M7").Value = -15088.091
$ws.Range("N7").Value = -4891.4287
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H122").Value = 3798.8125
$ws.Range("I122").Value = 3399.25
$ws.Range("K122").Value = 10197.75
$ws.Range("M122").Value = -7747.75
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 11104.056
$ws.Range("I126").Value = 15200.091
$ws.Range("J126").Value = 4667.4287
$ws.Range("K126").Value = 45600.273
$ws.Range("L126").Value = 14002.2861
$ws.Range("M126").Value = -43130.273
$ws.Range("N126").Value = -18942.2861
$ws.Range("H136").Value = 4032.138
$ws.Range("I136").Value = 3215.75
$ws.Range("K136").Value = 9647.25
$ws.Range("M136").Value = -7097.25
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1274.9231
$ws.Range("I122").Value = 1274.9231
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3824.7693
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1374.7693
$ws.Range("N122").ClearContents()
$ws.Range("H133").Value = 50715
$ws.Range("J133").Value = 50715
$ws.Range("L133").Value = 50715
$ws.Range("N133").Value = -60835
$ws.Range("H136").Value = 2533.72
$ws.Range("I136").Value = 1831.875
$ws.Range("K136").Value = 5495.625
$ws.Range("M136").Value = -2945.625
$ws.Range("N136").ClearContents()
